$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A10/B10 need to pick up the same cell style family as B1 (font id 2, the
# "plain" font) but with a General number format instead of the date format
# -- copy B1's format over, then force the number format back to General.
$ws.Range("B1").Copy()
$ws.Range("A10:B10").PasteSpecial(-4122)
$ws.Range("A10:B10").NumberFormat = "general"

# Update the "Vaccine tiêm trong ngày" row: clarify the label and split the
# single daily-dose total into dose-1 + dose-2 components (as text, like the
# other "x+y" breakdown cells in this sheet).
$ws.Range("A10").Value = "Vaccine tiêm trong ngày (mũi 1 + mũi 2)"
$ws.Range("B10").Value = "131274+10541"
